$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric; force text format so Excel
# keeps them as strings instead of silently converting to numbers.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D14",
    "D16",
    "D17",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D49",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptocurrency snapshot values
$ws.Range("D2").Value = "41.456.23"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.430.07"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +1.61%  "
$ws.Range("D5").Value = "308.20"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "88.53"
$ws.Range("E6").Value = "  -7.30%  "
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  -4.73%  "
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  -6.58%  "
$ws.Range("D10").Value = "31.50"
$ws.Range("E10").Value = "  -7.97%  "
$ws.Range("D11").Value = "0.0757"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "2.796.35"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  -5.51%  "
$ws.Range("D15").Value = "2.403.90"
$ws.Range("E15").Value = "  -4.58%  "
$ws.Range("D16").Value = "14.82"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "0.747"
$ws.Range("E17").Value = "  -5.63%  "
$ws.Range("D18").Value = "41.130.58"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "6.11"
$ws.Range("E19").Value = "  -4.56%  "
$ws.Range("D20").Value = "0.0₃0893"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").Value = "68.23"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").Value = "10.60"
$ws.Range("E22").Value = "  -9.80%  "
$ws.Range("D23").Value = "229.44"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").Value = "2.64"
$ws.Range("E24").Value = "  -5.68%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -6.31%  "
$ws.Range("D27").Value = "23.26"
$ws.Range("E27").Value = "  -6.07%  "
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").Value = "9.36"
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("D30").Value = "34.29"
$ws.Range("E30").Value = "  -7.30%  "
$ws.Range("D31").Value = "151.06"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").Value = "5.21"
$ws.Range("E32").Value = "  -8.17%  "
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").Value = "2.48"
$ws.Range("E34").Value = "  -5.59%  "
$ws.Range("D35").Value = "0.0733"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("D36").Value = "16.88"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").Value = "2.86"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("D38").Value = "1.75"
$ws.Range("E38").Value = "  -7.25%  "
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").Value = "0.0969"
$ws.Range("E40").Value = "  -8.28%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "3.90"
$ws.Range("E41").Value = "  -4.62%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.02"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").Value = "18.69"
$ws.Range("E43").Value = "  -13.40%  "
$ws.Range("D44").Value = "1.909.17"
$ws.Range("E44").Value = "  -4.88%  "
$ws.Range("D45").Value = "0.0273"
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  -9.32%  "
$ws.Range("D47").Value = "2.680.90"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "92.71"
$ws.Range("E49").Value = "  -5.37%  "
$ws.Range("E50").Value = "  -7.26%  "
$ws.Range("D51").Value = "70.64"
$ws.Range("E51").Value = "  -9.04%  "

# Restore default (General/Normal) formatting on the forced-text cells
# so no stray style/number-format differences are introduced.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}

